# Apply the authoring changes described in the commit:
# "Designed the basic UI of the chat bot on the web app."
#
# 1. Slide 3 ("First and foremost... What sets CryptoShare apart?"):
#    mention cross-platform functionality alongside open-source nature.
# 2. Slide 6 ("Completed User Stories"): add a new lead-in bullet
#    "Currently, users can..." (no bullet glyph, flush left) and let the
#    body shrink text on overflow (normAutofit).
# 3. Slide 8 ("Remaining User Stories"): rename "Budget income." to
#    "Budgeting system."

$p = $ppt.ActivePresentation

# --- Slide 3: investors paragraph -----------------------------------
$slide3 = $p.Slides.Item(3)
$investorsShape = $slide3.Shapes.Item(2)
$investorsRange = $investorsShape.TextFrame.TextRange
$investorsPara = $investorsRange.Paragraphs(2)
# Clear to unrelated text first so PowerPoint doesn't try to keep a
# shared-prefix run alive -- this yields one clean run, matching how a
# fresh retype of the sentence behaves.
$investorsPara.Text = "zzz_placeholder_zzz"
$investorsPara.Text = "Investors would value the open-source nature, cross-platform functionality, the mixture of stocks and cryptos, the ability to import/export user data securely and easily (based on a large amount of user feedback from the previous iteration of this software). "

# --- Slide 6: add "Currently, users can..." lead-in bullet ----------
$slide6 = $p.Slides.Item(6)
$storiesShape = $slide6.Shapes.Item(2)
$storiesRange = $storiesShape.TextFrame.TextRange
[void]$storiesRange.Paragraphs(1).InsertBefore("Currently, users can…`r")

$newLead = $storiesRange.Paragraphs(1)
$newLead.ParagraphFormat.Bullet.Visible = $false
$leadRuler = $storiesShape.TextFrame.Ruler.Levels.Item(1)
$leadRuler.LeftMargin = 0
$leadRuler.FirstMargin = 0

# Shrink text on overflow now that there's an extra line.
$storiesShape.TextFrame.AutoSize = 2

# --- Slide 8: "Budget income." -> "Budgeting system." ----------------
$slide8 = $p.Slides.Item(8)
$remainingShape = $slide8.Shapes.Item(2)
$remainingRange = $remainingShape.TextFrame.TextRange
$budgetPara = $remainingRange.Paragraphs(2)
$budgetPara.Text = "zzz_placeholder_zzz"
$budgetPara.Text = "Budgeting system."
